$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2026-03-01 -> 2026-03-02, serial 46082 -> 46083) for every data row
# (rows 2 through 216).
$ws.Range("C2:C216").Value = 46083
